$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 3).Value = 43806
$ws.Cells.Item(2, 4).Value = 65448919
$ws.Cells.Item(3, 3).Value = 103642
$ws.Cells.Item(3, 4).Value = 157513003
$ws.Cells.Item(4, 3).Value = 35345
$ws.Cells.Item(4, 4).Value = 55690380
$ws.Cells.Item(5, 3).Value = 10186
$ws.Cells.Item(5, 4).Value = 16602300
$ws.Cells.Item(6, 3).Value = 2648
$ws.Cells.Item(6, 4).Value = 4542691
$ws.Cells.Item(7, 3).Value = 300
$ws.Cells.Item(7, 4).Value = 561517
$ws.Cells.Item(12, 3).Value = 46511
$ws.Cells.Item(12, 4).Value = 63518046
$ws.Cells.Item(13, 3).Value = 11072
$ws.Cells.Item(13, 4).Value = 16323310
$ws.Cells.Item(14, 3).Value = 28984
$ws.Cells.Item(14, 4).Value = 43177958
$ws.Cells.Item(15, 3).Value = 9184
$ws.Cells.Item(15, 4).Value = 13970920
$ws.Cells.Item(16, 3).Value = 2448
$ws.Cells.Item(16, 4).Value = 3790379
$ws.Cells.Item(17, 3).Value = 527
$ws.Cells.Item(17, 4).Value = 829171
$ws.Cells.Item(20, 3).Value = 11406
$ws.Cells.Item(20, 4).Value = 15072545
$ws.Cells.Item(21, 3).Value = 15082
$ws.Cells.Item(21, 4).Value = 21935752
$ws.Cells.Item(22, 3).Value = 35009
$ws.Cells.Item(22, 4).Value = 51717992
$ws.Cells.Item(23, 3).Value = 11289
$ws.Cells.Item(23, 4).Value = 17084285
$ws.Cells.Item(24, 3).Value = 3016
$ws.Cells.Item(24, 4).Value = 4645654
$ws.Cells.Item(25, 3).Value = 666
$ws.Cells.Item(25, 4).Value = 1049301
$ws.Cells.Item(26, 3).Value = 56
$ws.Cells.Item(26, 4).Value = 113569
$ws.Cells.Item(27, 3).Value = 12895
$ws.Cells.Item(27, 4).Value = 17102974
$ws.Cells.Item(28, 3).Value = 8882
$ws.Cells.Item(28, 4).Value = 13115438
$ws.Cells.Item(29, 3).Value = 25393
$ws.Cells.Item(29, 4).Value = 37942939
$ws.Cells.Item(30, 3).Value = 8687
$ws.Cells.Item(30, 4).Value = 13240982
$ws.Cells.Item(31, 3).Value = 2225
$ws.Cells.Item(31, 4).Value = 3412322
$ws.Cells.Item(32, 3).Value = 474
$ws.Cells.Item(32, 4).Value = 752430
$ws.Cells.Item(34, 3).Value = 9330
$ws.Cells.Item(34, 4).Value = 12303513
$ws.Cells.Item(35, 3).Value = 3864
$ws.Cells.Item(35, 4).Value = 5725986
$ws.Cells.Item(36, 3).Value = 9027
$ws.Cells.Item(36, 4).Value = 13609236
$ws.Cells.Item(37, 3).Value = 3535
$ws.Cells.Item(37, 4).Value = 5424432
$ws.Cells.Item(38, 3).Value = 905
$ws.Cells.Item(38, 4).Value = 1382216
$ws.Cells.Item(39, 3).Value = 193
$ws.Cells.Item(39, 4).Value = 312686
$ws.Cells.Item(41, 3).Value = 2890
$ws.Cells.Item(41, 4).Value = 3918725
$ws.Cells.Item(42, 3).Value = 19688
$ws.Cells.Item(42, 4).Value = 29051467
$ws.Cells.Item(43, 3).Value = 57053
$ws.Cells.Item(43, 4).Value = 84938228
$ws.Cells.Item(44, 3).Value = 20857
$ws.Cells.Item(44, 4).Value = 31479258
$ws.Cells.Item(45, 3).Value = 6291
$ws.Cells.Item(45, 4).Value = 9537339
$ws.Cells.Item(46, 3).Value = 1546
$ws.Cells.Item(46, 4).Value = 2428276
$ws.Cells.Item(47, 3).Value = 102
$ws.Cells.Item(47, 4).Value = 197461
$ws.Cells.Item(50, 3).Value = 19044
$ws.Cells.Item(50, 4).Value = 25304450
$ws.Cells.Item(51, 3).Value = 2464
$ws.Cells.Item(51, 4).Value = 3669324
$ws.Cells.Item(52, 3).Value = 8246
$ws.Cells.Item(52, 4).Value = 12394528
$ws.Cells.Item(53, 3).Value = 2771
$ws.Cells.Item(53, 4).Value = 4328590
$ws.Cells.Item(54, 3).Value = 878
$ws.Cells.Item(54, 4).Value = 1367798
$ws.Cells.Item(55, 3).Value = 246
$ws.Cells.Item(55, 4).Value = 414110
$ws.Cells.Item(56, 3).Value = 27
$ws.Cells.Item(56, 4).Value = 66000
$ws.Cells.Item(57, 3).Value = 8161
$ws.Cells.Item(57, 4).Value = 11297552
$ws.Cells.Item(58, 3).Value = 1708
$ws.Cells.Item(58, 4).Value = 3512193
$ws.Cells.Item(59, 3).Value = 4061
$ws.Cells.Item(59, 4).Value = 8294542
$ws.Cells.Item(60, 3).Value = 1615
$ws.Cells.Item(60, 4).Value = 3324997
$ws.Cells.Item(61, 3).Value = 538
$ws.Cells.Item(61, 4).Value = 1092845
$ws.Cells.Item(64, 3).Value = 2650
$ws.Cells.Item(64, 4).Value = 5023301
$ws.Cells.Item(65, 3).Value = 17938
$ws.Cells.Item(65, 4).Value = 26748775
$ws.Cells.Item(66, 3).Value = 50762
$ws.Cells.Item(66, 4).Value = 76866629
$ws.Cells.Item(67, 3).Value = 17687
$ws.Cells.Item(67, 4).Value = 27516181
$ws.Cells.Item(68, 3).Value = 5192
$ws.Cells.Item(68, 4).Value = 8195786
$ws.Cells.Item(69, 3).Value = 1200
$ws.Cells.Item(69, 4).Value = 2059572
$ws.Cells.Item(70, 3).Value = 122
$ws.Cells.Item(70, 4).Value = 234582
$ws.Cells.Item(71, 3).Value = 19
$ws.Cells.Item(71, 4).Value = 25619
$ws.Cells.Item(73, 3).Value = 16842
$ws.Cells.Item(73, 4).Value = 22185213
$ws.Cells.Item(74, 3).Value = 67739
$ws.Cells.Item(74, 4).Value = 105999286
$ws.Cells.Item(75, 3).Value = 182072
$ws.Cells.Item(75, 4).Value = 288136551
$ws.Cells.Item(76, 3).Value = 77710
$ws.Cells.Item(76, 4).Value = 128134555
$ws.Cells.Item(77, 3).Value = 25726
$ws.Cells.Item(77, 4).Value = 44331916
$ws.Cells.Item(78, 3).Value = 7210
$ws.Cells.Item(78, 4).Value = 14163195
$ws.Cells.Item(79, 3).Value = 611
$ws.Cells.Item(79, 4).Value = 1609974
$ws.Cells.Item(85, 3).Value = 66330
$ws.Cells.Item(85, 4).Value = 91865385
$ws.Cells.Item(86, 3).Value = 5304
$ws.Cells.Item(86, 4).Value = 7751296
$ws.Cells.Item(87, 3).Value = 12948
$ws.Cells.Item(87, 4).Value = 19281008
$ws.Cells.Item(88, 3).Value = 4210
$ws.Cells.Item(88, 4).Value = 6357598
$ws.Cells.Item(89, 3).Value = 1487
$ws.Cells.Item(89, 4).Value = 2221611
$ws.Cells.Item(90, 3).Value = 366
$ws.Cells.Item(90, 4).Value = 570512
$ws.Cells.Item(93, 3).Value = 6026
$ws.Cells.Item(93, 4).Value = 8099791
$ws.Cells.Item(94, 3).Value = 1920
$ws.Cells.Item(94, 4).Value = 2820495
$ws.Cells.Item(95, 3).Value = 6109
$ws.Cells.Item(95, 4).Value = 9230878
$ws.Cells.Item(96, 3).Value = 2176
$ws.Cells.Item(96, 4).Value = 3324357
$ws.Cells.Item(97, 3).Value = 798
$ws.Cells.Item(97, 4).Value = 1217957
$ws.Cells.Item(98, 3).Value = 247
$ws.Cells.Item(98, 4).Value = 399597
$ws.Cells.Item(101, 3).Value = 4139
$ws.Cells.Item(101, 4).Value = 5533699
$ws.Cells.Item(102, 3).Value = 948
$ws.Cells.Item(102, 4).Value = 1863288
$ws.Cells.Item(103, 3).Value = 647
$ws.Cells.Item(103, 4).Value = 1356212
$ws.Cells.Item(107, 3).Value = 6
$ws.Cells.Item(107, 4).Value = 20490
$ws.Cells.Item(108, 3).Value = 12598
$ws.Cells.Item(108, 4).Value = 18893194
$ws.Cells.Item(109, 3).Value = 32441
$ws.Cells.Item(109, 4).Value = 48775678
$ws.Cells.Item(110, 3).Value = 10864
$ws.Cells.Item(110, 4).Value = 16605028
$ws.Cells.Item(111, 3).Value = 3068
$ws.Cells.Item(111, 4).Value = 4754131
$ws.Cells.Item(112, 3).Value = 634
$ws.Cells.Item(112, 4).Value = 1010851
$ws.Cells.Item(113, 3).Value = 94
$ws.Cells.Item(113, 4).Value = 213219
$ws.Cells.Item(116, 3).Value = 10856
$ws.Cells.Item(116, 4).Value = 14344366
$ws.Cells.Item(117, 3).Value = 34995
$ws.Cells.Item(117, 4).Value = 51434996
$ws.Cells.Item(118, 3).Value = 73846
$ws.Cells.Item(118, 4).Value = 109712185
$ws.Cells.Item(119, 3).Value = 23661
$ws.Cells.Item(119, 4).Value = 35761294
$ws.Cells.Item(120, 3).Value = 6780
$ws.Cells.Item(120, 4).Value = 10356144
$ws.Cells.Item(121, 3).Value = 1425
$ws.Cells.Item(121, 4).Value = 2292344
$ws.Cells.Item(122, 3).Value = 150
$ws.Cells.Item(122, 4).Value = 257291
$ws.Cells.Item(126, 3).Value = 28621
$ws.Cells.Item(126, 4).Value = 38241286
$ws.Cells.Item(127, 3).Value = 42077
$ws.Cells.Item(127, 4).Value = 62815435
$ws.Cells.Item(128, 3).Value = 87266
$ws.Cells.Item(128, 4).Value = 131993954
$ws.Cells.Item(129, 3).Value = 26716
$ws.Cells.Item(129, 4).Value = 41975066
$ws.Cells.Item(130, 3).Value = 7330
$ws.Cells.Item(130, 4).Value = 11628289
$ws.Cells.Item(131, 3).Value = 1644
$ws.Cells.Item(131, 4).Value = 2825647
$ws.Cells.Item(135, 3).Value = 35543
$ws.Cells.Item(135, 4).Value = 47390342
$ws.Cells.Item(136, 3).Value = 15309
$ws.Cells.Item(136, 4).Value = 22441150
$ws.Cells.Item(137, 3).Value = 36189
$ws.Cells.Item(137, 4).Value = 53550496
$ws.Cells.Item(138, 3).Value = 12781
$ws.Cells.Item(138, 4).Value = 19124406
$ws.Cells.Item(139, 3).Value = 3395
$ws.Cells.Item(139, 4).Value = 5176733
$ws.Cells.Item(140, 3).Value = 652
$ws.Cells.Item(140, 4).Value = 1056476
$ws.Cells.Item(141, 3).Value = 59
$ws.Cells.Item(141, 4).Value = 117989
$ws.Cells.Item(144, 3).Value = 12047
$ws.Cells.Item(144, 4).Value = 16099912
$ws.Cells.Item(145, 3).Value = 41617
$ws.Cells.Item(145, 4).Value = 63088046
$ws.Cells.Item(146, 3).Value = 95118
$ws.Cells.Item(146, 4).Value = 146889600
$ws.Cells.Item(147, 3).Value = 28461
$ws.Cells.Item(147, 4).Value = 45512580
$ws.Cells.Item(148, 3).Value = 7695
$ws.Cells.Item(148, 4).Value = 12767855
$ws.Cells.Item(149, 3).Value = 1991
$ws.Cells.Item(149, 4).Value = 3574375
$ws.Cells.Item(150, 3).Value = 177
$ws.Cells.Item(150, 4).Value = 402276
$ws.Cells.Item(151, 3).Value = 23
$ws.Cells.Item(151, 4).Value = 51500
$ws.Cells.Item(152, 3).Value = 33130
$ws.Cells.Item(152, 4).Value = 45027918
